$wb = $excel.ActiveWorkbook

# Update status text from "Ready for handoff" to "In Translation" on every sheet
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"

# The new status text is shorter, so the status columns are narrowed to fit
$wsOverview.Columns("E:F").ColumnWidth = 12.5
$wsZhCn.Columns("C:C").ColumnWidth = 12.5
$wsDeDe.Columns("C:C").ColumnWidth = 12.5
